# ADDENDUM_TEMPLATE_v2.docx edit: "combined deferred and compensation"
#
# 1. Shrink the "Date of Birth" row height (participant table) from
#    395 -> 260 twentieths of a point (19.75pt -> 13pt).
# 2. Split the *second* "alternate_payee.phone" run (the Phone row's
#    area-code placeholder, {{ phone_number_part(alternate_payee.phone, 1) }})
#    into "alternate_" + "payee.phone", relocating the lone "_GoBack"
#    bookmark to sit between the two halves.
# 3. Because a bookmark name is unique, adding "_GoBack" at its new spot
#    removes it from its old location (the empty paragraph near the end
#    of the document) automatically.

$d = $word.ActiveDocument

# --- 1. Row height -------------------------------------------------------
$participantTable = $d.Tables.Item(1)
for ($i = 1; $i -le $participantTable.Rows.Count; $i++) {
    $row = $participantTable.Rows.Item($i)
    if ($row.Cells.Item(1).Range.Text -like "Date of Birth:*") {
        $row.Height = 13.0
        break
    }
}

# --- 2/3. Split the 2nd "alternate_payee.phone" run & move the bookmark --
$alternateTable = $d.Tables.Item(2)
for ($i = 1; $i -le $alternateTable.Rows.Count; $i++) {
    $row = $alternateTable.Rows.Item($i)
    if ($row.Cells.Item(1).Range.Text -like "Phone:*") {
        $phoneCellRange = $row.Cells.Item(2).Range
        break
    }
}

$searchRange = $d.Range($phoneCellRange.Start, $phoneCellRange.End)

# Advance past the 1st occurrence, land the range on the 2nd occurrence.
$searchRange.Find.Execute("alternate_payee.phone", $true, $false, $false, `
    $false, $false, $true, 1, $false, "", 0) | Out-Null
$searchRange.Find.Execute("alternate_payee.phone", $true, $false, $false, `
    $false, $false, $true, 1, $false, "", 0) | Out-Null

# Split point: right after "alternate_" (10 characters in).
$splitPoint = $searchRange.Start + 10

# Move the existing "_GoBack" bookmark here (removes it from its old spot).
$bmRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# Rewrite the trailing half as its own run (drops the inherited rsid,
# matching the target markup) while keeping the same visible text. Text is
# routed through a throwaway placeholder first so the final assignment is
# never a same-text no-op (which the engine would otherwise skip).
$afterRange = $d.Range($splitPoint, $searchRange.End)
$afterRange.Text = "ZZPLACEHOLDERZZ"
$afterRange2 = $d.Range($splitPoint, $afterRange.End)
$afterRange2.Text = "payee.phone"
